# "Generate Report for Archive"
# Update the localization status from "Ready for handoff" to "In Translation"
# on every sheet that surfaces it, and shrink the now-narrower status
# columns to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: status cells E2 (zh-cn) and F2 (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Columns E and F narrow now that the status text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: status cell C2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: status cell C2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
